$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells before assignment to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '34.971.98'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.841.87'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '231.50'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '39.64'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '0.330'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("D11").Value = '0.0982'
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").Value = '2.109.39'
$ws.Range("E12").Value = '  +1.98%  '
$ws.Range("D13").Value = '11.46'
$ws.Range("E13").Value = '  +4.28%  '
$ws.Range("D14").Value = '1.842.66'
$ws.Range("E14").Value = '  +1.46%  '
$ws.Range("D15").Value = '0.672'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = '4.63'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '34.936.15'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '69.78'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").Value = '239.43'
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = '12.11'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").Value = '4.65'
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("D25").Value = '171.71'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '7.78'
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").Value = '17.42'
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("E28").Value = '  +3.38%  '
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = '0.0551'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").Value = '3.94'
$ws.Range("E32").Value = '  -3.11%  '
$ws.Range("D33").Value = '3.96'
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  +9.81%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +21.43%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '1.22'
$ws.Range("E36").Value = '  +7.46%  '
$ws.Range("D37").Value = '0.705'
$ws.Range("E37").Value = '  +2.56%  '
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").Value = '  +7.85%  '
$ws.Range("D39").Value = '89.91'
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("D40").Value = '1.350.95'
$ws.Range("E40").Value = '  +3.08%  '
$ws.Range("D41").Value = '0.0194'
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("D42").Value = '14.69'
$ws.Range("E42").Value = '  +2.36%  '
$ws.Range("D43").Value = '2.28'
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("D44").Value = '2.42'
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("E45").Value = '  +0.35%  '
$ws.Range("D46").Value = '0.0530'
$ws.Range("E46").Value = '  +3.31%  '
$ws.Range("D47").Value = '6.23'
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("D48").Value = '2.025.27'
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0678'
$ws.Range("E49").Value = '  +2.76%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '3.38'
$ws.Range("E51").Value = '  +20.54%  '

# Restore default style on Price cells (keeps them as plain text, no custom format)
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"

